$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update "想去人数" (F column) for rows 4 and 5
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1110
$wsExhibit.Range("F5").Value = 579

# Sheet "全部类型" (All Types): update "想去人数" (F column) for rows 4 and 6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1110
$wsAll.Range("F6").Value = 579
